# The struct_type label columns (B, D, F, H, J) on the data rows were
# imported with a stray ".1" suffix tacked onto the Chinese struct-type
# name (an artifact of the upstream export). Strip that trailing ".1"
# from each of those cells, row by row, leaving the numeric coeff
# columns (C, E, G, I, K) and everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$labelCols = @(2, 4, 6, 8, 10)  # B, D, F, H, J

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    foreach ($c in $labelCols) {
        $cell = $ws.Cells.Item($r, $c)
        $s = $cell.Value2
        if ($s -ne $null -and $s.EndsWith(".1")) {
            $cell.Value = $s.Substring(0, $s.Length - 2)
        }
    }
}
